# Apply the edits described by the diff for test_data/sample_data.xlsx:
#   1. Give column A an explicit custom width (was default/unset, now ~26.8 chars).
#   2. Flip the sign of C13 (profit figure) to a negative value.
#   3. Move the worksheet's active cell/selection from I12 to D14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Custom width for column A.
#    The target stored OOXML <col width="26.81640625"/> corresponds to the
#    Excel COM "ColumnWidth" (character-unit) value of 26.81640625 - 5/7
#    (Excel's column width model adds a fixed 5px/MaxDigitWidth padding term
#    on top of the character width when persisting to the file format).
$ws.Columns.Item(1).ColumnWidth = 26.81640625 - 5/7

# 2. C13 profit becomes negative.
$ws.Range("C13").Value = -7300

# 3. Update the active selection to D14.
$ws.Range("D14").Select()
